# Adds two new benchmark worksheets, "same_elements" and "partly_same",
# at the end of the workbook, using the same layout as the existing
# sheets: header row (array sizes) in B1:G1, then one row per data type
# (byte/int/string/date) in A2:G5 holding timing results. The newly
# added "partly_same" sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsSame = $wb.Worksheets.Add($null, $lastSheet)
$wsSame.Name = "same_elements"

$wsPartlySame = $wb.Worksheets.Add($null, $wsSame)
$wsPartlySame.Name = "partly_same"

# Column headers (array sizes) shared by every benchmark sheet.
$headers = @(5, 50, 500, 5000, 50000, 500000)

# Row labels (data types under test).
$labels = @("byte", "int", "string", "date")

$sameElementsData = @(
    @(0, 0.00050000000000000001, 0,                      0.00050100000000000003, 0, 0),
    @(0, 0,                      0,                      0.0010009999999999999,  0, 0),
    @(0, 0,                      0.00050000000000000001, 0.00050000000000000001, 0, 0),
    @(0, 0,                      0,                      0.0010009999999999999,  0, 0)
)

$partlySameData = @(
    @(0, 0.00050100000000000003, 0.028561,             3.3049529999999998, 0, 0),
    @(0, 0.00050000000000000001, 0.031033000000000002, 3.6903679999999999, 0, 0),
    @(0, 0,                      0.032534,             3.8209919999999999, 0, 0),
    @(0, 0.00050199999999999995, 0.030568000000000001, 3.6833469999999999, 0, 0)
)

$sheetConfigs = @(
    @{ Sheet = $wsSame; Data = $sameElementsData },
    @{ Sheet = $wsPartlySame; Data = $partlySameData }
)

foreach ($cfg in $sheetConfigs) {
    $ws = $cfg.Sheet
    $rows = $cfg.Data

    for ($c = 0; $c -lt $headers.Length; $c++) {
        $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
    }

    for ($r = 0; $r -lt $labels.Length; $r++) {
        $ws.Cells.Item($r + 2, 1).Value = $labels[$r]
        $rowVals = $rows[$r]
        for ($c = 0; $c -lt $rowVals.Length; $c++) {
            $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
        }
    }
}

# The newly added "partly_same" sheet becomes the active/selected tab,
# with the selection left wherever the user last clicked (L13).
$wsPartlySame.Activate()
$wsPartlySame.Range("L13").Select()
